# Tasks solved at 08.02.2022
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the actual start date for week 2 (row 3) in column C,
# matching the date-formatted style already used in column B / C2.
$ws.Range("C3").Value = 44600
$ws.Range("C3").NumberFormat = "m/d/yy"

# Move the selection to D3 (matches the saved selection state)
$ws.Range("D3").Select()
